$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12, shifting the "R1 to ISP Router" / "R2 to ISP
# Router" rows down to 13/14.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the Radius Server details (written in the
# same order as the original edit so shared-string indices line up).
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Radius Server"
$ws.Range("C12").Value = "172.16.8.224"
$ws.Range("F12").Value = "172.16.8.225"
$ws.Range("D12").Value = "/30"
$ws.Range("E12").Value = "255.255.255.252"
$ws.Range("G12").Value = "172.16.8.226"
$ws.Range("H12").Value = "172.16.8.227"
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 2

$ws.Range("F15").Select()
